$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2649.9756
$ws.Range("I62").Value = 2579.276
$ws.Range("J62").Value = 2820.8333
$ws.Range("K62").Value = 2579.276
$ws.Range("L62").Value = 2820.8333
$ws.Range("M62").Value = -1955.276
$ws.Range("N62").Value = -4068.8333

$ws.Range("H65").Value = 2649.9756
$ws.Range("I65").Value = 2579.276
$ws.Range("J65").Value = 2820.8333
$ws.Range("K65").Value = 12896.38
$ws.Range("L65").Value = 14104.1665
$ws.Range("M65").Value = -9776.379999999999
$ws.Range("N65").Value = -20344.1665

$ws.Range("H125").Value = 1849.0625
$ws.Range("I125").Value = 1259
$ws.Range("J125").Value = 2117.2727
$ws.Range("K125").Value = 11331
$ws.Range("L125").Value = 19055.4543
$ws.Range("M125").Value = -8871
$ws.Range("N125").Value = -23975.4543

$ws.Range("H137").Value = 1082.3438
$ws.Range("I137").Value = 963.7143
$ws.Range("J137").Value = 1174.6111
$ws.Range("K137").Value = 2891.1429
$ws.Range("L137").Value = 3523.8333
$ws.Range("M137").Value = -341.1428999999998
$ws.Range("N137").Value = -8623.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 50000000
$ws.Range("I8").Value = 50000000
$ws.Range("K8").Value = 50000000
$ws.Range("M8").Value = -49999856

$ws.Range("H32").Value = 515128.97
$ws.Range("I32").Value = 4106.8696
$ws.Range("J32").Value = 2865830.8
$ws.Range("K32").Value = 4106.8696
$ws.Range("L32").Value = 2865830.8
$ws.Range("M32").Value = -3819.8696
$ws.Range("N32").Value = -2866404.8

$ws.Range("H74").Value = 1163.7188
$ws.Range("I74").Value = 1188.2413
$ws.Range("J74").Value = 926.6667
$ws.Range("K74").Value = 1188.2413
$ws.Range("L74").Value = 926.6667
$ws.Range("M74").Value = -314.2412999999999
$ws.Range("N74").Value = -2674.6667

$ws.Range("H77").Value = 1163.7188
$ws.Range("I77").Value = 1188.2413
$ws.Range("J77").Value = 926.6667
$ws.Range("K77").Value = 5941.206499999999
$ws.Range("L77").Value = 4633.3335
$ws.Range("M77").Value = -1573.206499999999
$ws.Range("N77").Value = -13369.3335

$ws.Range("H132").Value = 1953.0682
$ws.Range("I132").Value = 1159.2593
$ws.Range("J132").Value = 3213.8235
$ws.Range("K132").Value = 3477.7779
$ws.Range("L132").Value = 9641.470499999999
$ws.Range("M132").Value = -947.7779
$ws.Range("N132").Value = -14701.4705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11439.5
$ws.Range("I82").Value = 3426.625
$ws.Range("J82").Value = 22123.334
$ws.Range("K82").Value = 3426.625
$ws.Range("L82").Value = 22123.334
$ws.Range("M82").Value = -3043.625
$ws.Range("N82").Value = -22889.334

$ws.Range("H85").Value = 11439.5
$ws.Range("I85").Value = 3426.625
$ws.Range("J85").Value = 22123.334
$ws.Range("K85").Value = 3426.625
$ws.Range("L85").Value = 22123.334
$ws.Range("M85").Value = -2100.625
$ws.Range("N85").Value = -24775.334

$ws.Range("H134").Value = 1511.7539
$ws.Range("I134").Value = 1184.3556
$ws.Range("J134").Value = 2248.4
$ws.Range("K134").Value = 3553.066800000001
$ws.Range("L134").Value = 6745.200000000001
$ws.Range("M134").Value = -1018.066800000001
$ws.Range("N134").Value = -11815.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 302.72
$ws.Range("I22").Value = 260.85715
$ws.Range("J22").Value = 356
$ws.Range("K22").Value = 260.85715
$ws.Range("L22").Value = 356
$ws.Range("M22").Value = 89.14285000000001
$ws.Range("N22").Value = -1056

$ws.Range("H80").Value = 26533.334
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 26533.334
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H132").Value = 1884.4147
$ws.Range("I132").Value = 1565.5358
$ws.Range("J132").Value = 2571.2307
$ws.Range("K132").Value = 4696.607400000001
$ws.Range("L132").Value = 7713.6921
$ws.Range("M132").Value = -2166.607400000001
$ws.Range("N132").Value = -12773.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 106
$ws.Range("I14").Value = 106
$ws.Range("K14").Value = 318
$ws.Range("M14").Value = -145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 200800000
$ws.Range("J33").Value = 200800000
$ws.Range("L33").Value = 200800000
$ws.Range("N33").Value = -200800504

$ws.Range("H70").Value = 8279501
$ws.Range("I70").Value = 18154050
$ws.Range("J70").Value = 6229.2705
$ws.Range("K70").Value = 18154050
$ws.Range("L70").Value = 6229.2705
$ws.Range("M70").Value = -18153780
$ws.Range("N70").Value = -6769.2705

$ws.Range("H73").Value = 8279501
$ws.Range("I73").Value = 18154050
$ws.Range("J73").Value = 6229.2705
$ws.Range("K73").Value = 18154050
$ws.Range("L73").Value = 6229.2705
$ws.Range("M73").Value = -18153114
$ws.Range("N73").Value = -8101.2705

$ws.Range("H132").Value = 2186.75
$ws.Range("I132").Value = 2001.3043
$ws.Range("J132").Value = 3039.8
$ws.Range("K132").Value = 6003.9129
$ws.Range("L132").Value = 9119.400000000001
$ws.Range("M132").Value = -3473.9129
$ws.Range("N132").Value = -14179.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2416.4211
$ws.Range("I61").Value = 1672.909
$ws.Range("J61").Value = 3438.75
$ws.Range("K61").Value = 1672.909
$ws.Range("L61").Value = 3438.75
$ws.Range("M61").Value = -1470.909
$ws.Range("N61").Value = -3842.75

$ws.Range("H100").Value = 2205.2856
$ws.Range("I100").Value = 1383.25
$ws.Range("K100").Value = 1383.25
$ws.Range("M100").Value = -842.25

$ws.Range("H113").Value = 2416.4211
$ws.Range("I113").Value = 1672.909
$ws.Range("J113").Value = 3438.75
$ws.Range("K113").Value = 1672.909
$ws.Range("L113").Value = 3438.75
$ws.Range("M113").Value = 497.0909999999999
$ws.Range("N113").Value = -7778.75

$ws.Range("H132").Value = 2394.1606
$ws.Range("I132").Value = 2243.2173
$ws.Range("J132").Value = 3088.5
$ws.Range("K132").Value = 6729.651899999999
$ws.Range("L132").Value = 9265.5
$ws.Range("M132").Value = -4199.651899999999
$ws.Range("N132").Value = -14325.5

$ws.Range("H136").Value = 3246.0303
$ws.Range("I136").Value = 2053.1365
$ws.Range("J136").Value = 5631.8184
$ws.Range("K136").Value = 6159.4095
$ws.Range("L136").Value = 16895.4552
$ws.Range("M136").Value = -3609.4095
$ws.Range("N136").Value = -21995.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41668324
$ws.Range("I132").Value = 48388540
$ws.Range("J132").Value = 2989.6
$ws.Range("K132").Value = 145165620
$ws.Range("L132").Value = 8968.799999999999
$ws.Range("M132").Value = -145163090
$ws.Range("N132").Value = -14028.8
